# Contract creation workflow update:
#  - Row 2 (pedido 4503273185): COND.PAG moves from "ZCTP" to "Z000";
#    NV CONTRATO gets the newly created contract number; ITEM CONT. and
#    NV PEDIDO are cleared (no longer carried on this line).
#  - Row 3 (pedido 4503274856): COND.PAG moves from "A003" to "Z000";
#    TP CONTRATO changes to "ZDDR"; ITEM CONT., NV CONTRATO and NV PEDIDO
#    are cleared (contract not yet created for this line).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = "Z000"
$ws.Range("N2").ClearContents()
$ws.Range("P2").Value = 4600244273
$ws.Range("Q2").ClearContents()

# Row 3
$ws.Range("E3").Value = "Z000"
$ws.Range("M3").Value = "ZDDR"
$ws.Range("N3").ClearContents()
$ws.Range("P3").ClearContents()
$ws.Range("Q3").ClearContents()

# Restore the default top-left view and move the active selection, as left
# by the author after finishing the edits.
$ws.Range("N7").Select()
